$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: replace "minMerge" function row with "removeContradiction" row,
#     filling in the previously-empty Description / Unittest / PythonFile cells.
#     D22/F22 pick up the "filled-in description row" look (same formatting as
#     the analogous D4/F4 cells), so copy that formatting across first. ---
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C22").Value = "removeContradiction"
$ws.Range("D22").Value = 'Considers the independent rule (consisting of a single inequality or a sympy Or relational containing multiple inequalities) being proposed by one or more disciplines and merges them such that from the top-level, the rule(s) do not contradict each other'
$ws.Range("E22").Value = "Yes"
$ws.Range("F22").Value = "Necessary"
$ws.Rows(22).RowHeight = 57.6

# --- Row 25: fill in the previously-empty Description / Unittest / PythonFile cells
#     for the "AdjustCriteria" function row, matching the formatting used by the
#     analogous D5/E5/F5 cells. ---
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null
$ws.Range("F5").Copy() | Out-Null
$ws.Range("F25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("D25").Value = 'Cycles through the established criteria for allowing a space reduction that specifically pertain to the area of a discipline''s design space being reduced and relaxes one criterion when a space reduction is being forced for the discipline'
$ws.Range("E25").Value = "Yes"
$ws.Range("F25").Value = "Necessary"
$ws.Rows(25).RowHeight = 43.2

# --- Remove the now-vacant row 26 (it never held data), shifting the
#     EXPLORATION section and everything below it up by one row ---
$ws.Rows("26:26").Delete() | Out-Null

# --- Widen column C slightly to fit the new "removeContradiction" text ---
$ws.Columns(3).ColumnWidth = 16.833333333333336

# --- Restore the selection cell shown in the (still frozen-at-row-1) pane ---
$ws.Range("G22").Select() | Out-Null
